$d = $word.ActiveDocument
$d.Content.Find.Execute("Apéndice A: Cálculo de parámetros del modelo dinámico", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Apéndice A: Cálculo de parámetros de inercia del modelo dinámico", 2)
